$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6 and 7: relabel the "*_server" parameter names so they match the
# plain "host"/"user" labels used in rows 2/3 (host_server -> host, user_server -> user)
$ws.Range("A6").Value = "host"
$ws.Range("A7").Value = "user"

# Row 5: drop the trailing slash from the workdirectory path for the "Usar=True" block
$ws.Range("B5").Value = "C:/Users/iachenbach/Gobierno de la Ciudad de Buenos Aires/Pablo Alfredo Gadea - Tablero Facoep P BI/FACOEP/DBA/Reportes BI/2021/Monitoreo CRGs"

# Widen column B to fit the longer path text
$ws.Columns.Item(2).ColumnWidth = 134.8

# Move the active selection from B14 to A7
$ws.Range("A7").Select()
